# "Hjemme passive tweaks lichtwark deleted values"
# Update the B1:E1 header values, drop the stale/duplicate Lichtwark
# "CON" row (B2:E2) and the leading "STR" values (B3), and refresh the
# remaining STR row values (C3:E3) with the recomputed figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 (header counts) tweaks
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2 ("CON"): delete the now-invalid leading values entirely
$ws.Range("B2:E2").ClearContents()

# Row 3 ("STR"): drop B3, recompute/replace C3:E3
$ws.Range("B3").ClearContents()
$ws.Range("C3").Value = 1.2961979451674353
$ws.Range("D3").Value = -5.1271008514283096
$ws.Range("E3").Value = 11.095753804291924

# Leave the selection matching the trimmed data block
$ws.Range("B1:E3").Select()
